$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto symbol price/volume values (GitHub Actions data refresh)
$updates = @{
    'D2' = '293.01'
    'E2' = '-2.99%'
    'D3' = '31.38'
    'E3' = '-1.69%'
    'D4' = '4.974'
    'E4' = '-0.87%'
    'D5' = '0.07382'
    'E5' = '-5.51%'
    'D6' = '1.807'
    'E6' = '-11.66%'
    'D7' = '7.659'
    'E7' = '-1.76%'
    'D8' = '3.749'
    'E8' = '-0.72%'
    'D9' = '0.9104'
    'E9' = '-0.93%'
    'D10' = '0.1648'
    'E10' = '-5.54%'
    'D11' = '0.07617'
    'E11' = '-3.12%'
    'D12' = '0.08127'
    'E12' = '-7.16%'
    'D13' = '0.02990'
    'E13' = '-3.46%'
    'D14' = '0.09966'
    'E14' = '-0.35%'
    'D15' = '0.001493'
    'E15' = '-1.98%'
    'D16' = '0.005682'
    'E16' = '-4.18%'
    'D18' = '3.469'
    'E18' = '0.25%'
    'E19' = '-6.31%'
    'D21' = '0.1314'
    'E21' = '1.79%'
    'D22' = '4.324'
    'E22' = '3.60%'
    'D24' = '0.04488'
    'E24' = '-2.35%'
    'D25' = '0.001224'
    'E25' = '-1.38%'
    'D26' = '0.004052'
    'E26' = '-9.39%'
    'E27' = '0.03%'
    'D39' = '0.01630'
    'E39' = '-6.20%'
    'D40' = '0.04385'
    'E40' = '-7.57%'
    'D41' = '0.007419'
    'E41' = '2.34%'
    'E42' = '-2.36%'
    'E43' = '-5.27%'
    'D44' = '0.01002'
    'E44' = '-6.83%'
    'E45' = '-1.21%'
    'E46' = '0.01%'
    'D47' = '1.891'
    'E47' = '60.83%'
    'D48' = '0.003001'
    'E48' = '-14.52%'
    'E49' = '0.01%'
    'E50' = '0.01%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
